# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Cells are plain text in the sheet (Price looks numeric but is stored as
# text, e.g. "585.02"), so each write forces a Text number format first and
# restores the default "Normal" style afterwards to avoid Excel silently
# re-typing the cell as a number (which would also lose trailing zeros and
# introduce float rounding) or leaving a stray quote-prefix style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "71.036.32"
Set-TextCell "E2" "  +2.52%  "
Set-TextCell "D3" "3.574.67"
Set-TextCell "E3" "  +1.96%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.21%  "
Set-TextCell "D5" "585.02"
Set-TextCell "E5" "  +1.68%  "
Set-TextCell "D6" "190.62"
Set-TextCell "E6" "  +2.74%  "
Set-TextCell "D7" "0.626"
Set-TextCell "E7" "  +2.20%  "
Set-TextCell "D8" "3.567.04"
Set-TextCell "E8" "  +2.00%  "
Set-TextCell "E9" "  -0.11%  "
Set-TextCell "D10" "0.220"
Set-TextCell "E10" "  +15.94%  "
Set-TextCell "D11" "0.653"
Set-TextCell "E11" "  +0.82%  "
Set-TextCell "D12" "54.80"
Set-TextCell "E12" "  +1.24%  "
Set-TextCell "D13" "0.0000319"
Set-TextCell "E13" "  +5.86%  "
Set-TextCell "D14" "9.56"
Set-TextCell "E14" "  +1.19%  "
Set-TextCell "D15" "4.142.60"
Set-TextCell "E15" "  +1.86%  "
Set-TextCell "D16" "71.014.59"
Set-TextCell "E16" "  +2.52%  "
Set-TextCell "D17" "12.87"
Set-TextCell "E17" "  +4.61%  "
Set-TextCell "D18" "19.23"
Set-TextCell "E18" "  -0.66%  "
Set-TextCell "D19" "3.543.86"
Set-TextCell "E19" "  +1.24%  "
Set-TextCell "D20" "567.84"
Set-TextCell "E20" "  +4.14%  "
Set-TextCell "E21" "  +0.68%  "
Set-TextCell "E22" "  -0.32%  "
Set-TextCell "E23" "  -4.32%  "
Set-TextCell "D24" "4.61"
Set-TextCell "E24" "  +4.04%  "
Set-TextCell "D25" "4.91"
Set-TextCell "E25" "  -1.08%  "
Set-TextCell "D26" "94.40"
Set-TextCell "E26" "  +0.36%  "
Set-TextCell "D27" "11.24"
Set-TextCell "E27" "  -0.51%  "
Set-TextCell "E28" "  -0.14%  "
Set-TextCell "E29" "  +2.15%  "
Set-TextCell "D30" "32.66"
Set-TextCell "E30" "  +2.52%  "
Set-TextCell "E31" "  -0.45%  "
Set-TextCell "D32" "12.36"
Set-TextCell "E32" "  -1.86%  "
Set-TextCell "E33" "  +2.63%  "
Set-TextCell "D34" "64.05"
Set-TextCell "E34" "  -0.75%  "
Set-TextCell "D35" "3.79"
Set-TextCell "E35" "  +24.10%  "
Set-TextCell "D36" "3.28"
Set-TextCell "E36" "  +6.56%  "
Set-TextCell "B37" "Bittensor"
Set-TextCell "C37" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D37" "534.72"
Set-TextCell "E37" "  -1.11%  "
Set-TextCell "B38" "TheGraph"
Set-TextCell "C38" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D38" "0.412"
Set-TextCell "E38" "  +2.52%  "
Set-TextCell "D39" "38.52"
Set-TextCell "E39" "  +1.54%  "
Set-TextCell "B40" "PEPE"
Set-TextCell "C40" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D40" "0.0₃0803"
Set-TextCell "E40" "  +4.61%  "
Set-TextCell "B41" "Maker"
Set-TextCell "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D41" "3.650.51"
Set-TextCell "E41" "  +10.30%  "
Set-TextCell "E42" "  -0.08%  "
Set-TextCell "E43" "  +4.42%  "
Set-TextCell "E44" "  +2.85%  "
Set-TextCell "D45" "0.0473"
Set-TextCell "E45" "  +6.25%  "
Set-TextCell "E46" "  -1.16%  "
Set-TextCell "D47" "3.45"
Set-TextCell "E47" "  -0.26%  "
Set-TextCell "D48" "9.32"
Set-TextCell "E48" "  +4.48%  "
Set-TextCell "E49" "  +3.40%  "
Set-TextCell "E50" "  +0.02%  "
Set-TextCell "E51" "  +7.96%  "
